$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data lives in rows 2-101 (row 1 is the header row: volume / jarak / pixel).
# Column B ("volume") is rescaled by /1,000,000 and column C ("jarak") by /100,
# while columns A and D are left untouched.
for ($r = 2; $r -le 101; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value2 = $bCell.Value2 / 1000000

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value2 = $cCell.Value2 / 100
}
